$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sponsored`nLödha Bhandup West | Luxury 2, 2.5 & 3 BHK`nl-bhandup.com`nhttps://www.l-bhandup.com › official`nLödha on LBS Rd an upcoming residential project with the best of amenities ...`nBrochure & Floor Plans · Price · View Project Highlights · Connect on Whatsapp"
$ws.Range("B2").Value = "Rioga Premium Real Estate Advisory LLP"

$ws.Range("A3").Value = "Sponsored`nLodha New Launch Bhandup West | Download E-Brochure`nproptigermumbai.com`nhttps://www.proptigermumbai.com`nits Big Apartments with 2/3 BHK, its prime Location along the LBS Road, Starts @ ₹ 2.29Cr*`nPricing & Floor Plan · Download Brochure · Platinum Group · Transcon Developers"

$ws.Range("A4").Value = "Sponsored`nNew Launch At Bhandup`nlódháhomz.site`nhttps://www.lódháhomz.site › site-visit › enquire-now`nNew Launch At Bhandup, Mumbai — Its Big Apartments With 2/3 BHK, Its Prime Location Along The LBS Road, Starts @ ₹ 2.29cr*`nE-Brochure · Floor Plan · Pricing / Costing · Project Overview"
$ws.Range("B4").Value = "DIGITAL RUBIX"

$ws.Range("A5").Value = "Sponsored`nLodha® LBS is Coming Soon | Near upcoming metro station`nlodhagroup.in`nhttps://www.lodhagroup.in`nLive an exceptional lifestyle with forest living as Lodha comes to the prime LBS Marg."
$ws.Range("B5").Value = "Macrotech Developers Limited"

$ws.Range("A6").Value = "Sponsored`nLodha Bhandup West Mumbai | 2&3 BHK Apartment @ ₹2.29 CR*`nindia-sales-office.com`nhttps://www.india-sales-office.com › official-site`nStarts at ₹2.29 CR* | The prime LBS Road, Mumbai | Modern Homes | Lots of Space & Comfort. New launch project 2&3 BHK luxury apartment by Lodha Group, location at Bhandup, Mumbai. Best Price Guaranteed.`nFloor Plans Details · View Price & Payment Plan · Virtual Tour & Video · Request Brochure`nCall us"
$ws.Range("B6").Value = "Property Junction International Real Estate Broker LLC"
$ws.Range("C6").Value = "United Arab Emirates"

$ws.Range("A7").Value = "Sponsored`nLodha Bhandup | Lodha in Bhandup West`nnewlaunch-property.net`nhttps://www.newlaunch-property.net › lodha_bhandup`nGet Huge Discount On Booking. Schedule Your Free Site Visit & Get Complete Project Details. Buy 2/3 BHK Premium Apartments by Lodha Property in Bhandup West, Mumbai @ ₹ 2.29 Cr*. Get Floor Plan. Download E Brochure.`nFloor Plans · Price List · Project Highlights`nCall us"
$ws.Range("B7").Value = "GTF Technologies"

$ws.Range("A8").Value = "Sponsored`nLodha Bhandup By Lodha group | 2 & 3 BHK Homes at ₹ 2.29 Cr*`nnewlaunch-bhandup.co.in`nhttps://www.newlaunch-bhandup.co.in › new_launch`nEstablished in 2012`n—`nLuxury 2 & 3 BHK Homes Starts From ₹ 2.29 Cr* All Incl at Lodha Bhandup West, Mumbai. The project has easy access to both the Eastern Express Highway and the Eastern Motorway. Get The Best Rate.`nSite & Floor Plans · Proposed Floor Plan · Pricing · Proposed Site Plan · View Gallery`nCall us"
$ws.Range("B8").Value = "Propertypistol Realty Pvt ltd"
